# 自动更新Excel文件 - 2025-10-17 23:11:29
# For every data row: one day elapses, so "剩余" (days remaining, column E)
# drops by 1. When the countdown would hit 0, the cycle restarts: 剩余
# resets back to 总天 (total days, column D) and 开始时间 (start date,
# column F) rolls forward by that many days (YYYYMMDD arithmetic).
# Rows whose start date isn't a clean 8-digit YYYYMMDD value are left
# untouched (can't safely date-shift them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($row = 2; $row -le $lastRow; $row++) {
    $total = $ws.Cells.Item($row, 4).Value2
    $remaining = $ws.Cells.Item($row, 5).Value2
    $startDate = $ws.Cells.Item($row, 6).Value2

    if ($remaining -eq $null -or $total -eq $null -or $startDate -eq $null) { continue }

    # Every row's start date gets validated/parsed up front (needed in case
    # this row's countdown is about to roll over); a malformed date means
    # the row can't be safely processed, so it's left untouched entirely.
    $dateText = [string]$startDate
    if ($dateText.Length -ne 8) {
        continue
    }

    $year = [int]$dateText.Substring(0, 4)
    $month = [int]$dateText.Substring(4, 2)
    $day = [int]$dateText.Substring(6, 2)
    $dt = Get-Date -Year $year -Month $month -Day $day

    $newRemaining = $remaining - 1

    if ($newRemaining -le 0) {
        $dt = $dt.AddDays($total)
        $newStartDate = [int]$dt.ToString("yyyyMMdd")

        $ws.Cells.Item($row, 5).Value2 = $total
        $ws.Cells.Item($row, 6).Value2 = $newStartDate
    }
    else {
        $ws.Cells.Item($row, 5).Value2 = $newRemaining
    }
}
